$d = $word.ActiveDocument
$content = $d.Content
$xml = $content.WordOpenXML

$old3 = '<w:r><w:t>Edit launchsetting.json to not open in browser</w:t></w:r>'
$new3 = '<w:r><w:t xml:space="preserve">Edit </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>launchsetting.json</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to not open in browser</w:t></w:r>'
if ($xml.IndexOf($old3) -lt 0) { throw 'pattern not found for para 3' }
$xml = $xml.Replace($old3, $new3)

$old6 = '<w:r w:rsidRPr="00337E6F"><w:t>AutoMapper.Extensions.Microsoft.DependencyInjection</w:t></w:r>'
$new6 = '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00337E6F"><w:t>AutoMapper.Extensions.Microsoft.DependencyInjection</w:t></w:r><w:proofErr w:type="spellEnd"/>'
if ($xml.IndexOf($old6) -lt 0) { throw 'pattern not found for para 6' }
$xml = $xml.Replace($old6, $new6)

$old7 = '<w:r w:rsidRPr="00337E6F"><w:t>Microsoft.AspNetCore.Authentication.JwtBearer</w:t></w:r>'
$new7 = '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00337E6F"><w:t>Microsoft.AspNetCore.Authentication.JwtBearer</w:t></w:r><w:proofErr w:type="spellEnd"/>'
if ($xml.IndexOf($old7) -lt 0) { throw 'pattern not found for para 7' }
$xml = $xml.Replace($old7, $new7)

$old8 = '<w:r w:rsidRPr="00337E6F"><w:t>Microsoft.AspNetCore.Identity.EntityFrameworkCore</w:t></w:r>'
$new8 = '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00337E6F"><w:t>Microsoft.AspNetCore.Identity.EntityFrameworkCore</w:t></w:r><w:proofErr w:type="spellEnd"/>'
if ($xml.IndexOf($old8) -lt 0) { throw 'pattern not found for para 8' }
$xml = $xml.Replace($old8, $new8)

$old9 = '<w:r w:rsidRPr="00337E6F"><w:t>Microsoft.EntityFrameworkCore</w:t></w:r>'
$new9 = '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00337E6F"><w:t>Microsoft.EntityFrameworkCore</w:t></w:r><w:proofErr w:type="spellEnd"/>'
if ($xml.IndexOf($old9) -lt 0) { throw 'pattern not found for para 9' }
$xml = $xml.Replace($old9, $new9)

$old10 = '<w:r w:rsidRPr="00337E6F"><w:t>Microsoft.EntityFrameworkCore.SqlServer</w:t></w:r>'
$new10 = '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00337E6F"><w:t>Microsoft.EntityFrameworkCore.SqlServer</w:t></w:r><w:proofErr w:type="spellEnd"/>'
if ($xml.IndexOf($old10) -lt 0) { throw 'pattern not found for para 10' }
$xml = $xml.Replace($old10, $new10)

$old11 = '<w:r w:rsidRPr="00337E6F"><w:t>Microsoft.EntityFrameworkCore.Tools</w:t></w:r>'
$new11 = '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00337E6F"><w:t>Microsoft.EntityFrameworkCore.Tools</w:t></w:r><w:proofErr w:type="spellEnd"/>'
if ($xml.IndexOf($old11) -lt 0) { throw 'pattern not found for para 11' }
$xml = $xml.Replace($old11, $new11)

$old13 = '<w:r><w:t>Create ApplicationDbContext class</w:t></w:r>'
$new13 = '<w:r><w:t xml:space="preserve">Create </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ApplicationDbContext</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class</w:t></w:r>'
if ($xml.IndexOf($old13) -lt 0) { throw 'pattern not found for para 13' }
$xml = $xml.Replace($old13, $new13)

$old14 = '<w:r><w:t>Add connections string to appsettings.json</w:t></w:r>'
$new14 = '<w:r><w:t xml:space="preserve">Add connections string to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>appsettings.json</w:t></w:r><w:proofErr w:type="spellEnd"/>'
if ($xml.IndexOf($old14) -lt 0) { throw 'pattern not found for para 14' }
$xml = $xml.Replace($old14, $new14)

$old17 = '<w:r><w:t>Create ServiceExtensions class in Extensions folder</w:t></w:r>'
$new17 = '<w:r><w:t xml:space="preserve">Create </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ServiceExtensions</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class in Extensions folder</w:t></w:r>'
if ($xml.IndexOf($old17) -lt 0) { throw 'pattern not found for para 17' }
$xml = $xml.Replace($old17, $new17)

$old18 = '<w:r><w:t>Add CofigureSqlContext method to ServiceExtensions class</w:t></w:r>'
$new18 = '<w:r><w:t xml:space="preserve">Add </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CofigureSqlContext</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ServiceExtensions</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class</w:t></w:r>'
if ($xml.IndexOf($old18) -lt 0) { throw 'pattern not found for para 18' }
$xml = $xml.Replace($old18, $new18)

$old19 = '<w:r><w:t>Register SQL configuration in the ConfigureServices method in the Startup class</w:t></w:r>'
$new19 = '<w:r><w:t xml:space="preserve">Register SQL configuration in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ConfigureServices</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method in the Startup class</w:t></w:r>'
if ($xml.IndexOf($old19) -lt 0) { throw 'pattern not found for para 19' }
$xml = $xml.Replace($old19, $new19)

$old22 = '<w:r><w:t>Add CORS configuration in the ServiceExtensions class</w:t></w:r>'
$new22 = '<w:r><w:t xml:space="preserve">Add CORS configuration in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ServiceExtensions</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class</w:t></w:r>'
if ($xml.IndexOf($old22) -lt 0) { throw 'pattern not found for para 22' }
$xml = $xml.Replace($old22, $new22)

$old26 = '<w:r><w:t>Create User model class that extends IdentityUser</w:t></w:r>'
$new26 = '<w:r><w:t xml:space="preserve">Create User model class that extends </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IdentityUser</w:t></w:r><w:proofErr w:type="spellEnd"/>'
if ($xml.IndexOf($old26) -lt 0) { throw 'pattern not found for para 26' }
$xml = $xml.Replace($old26, $new26)

$old27 = '<w:r><w:t>Update ApplicationDbContext IdentityDbContext parent to use User class as IdentityUser extender</w:t></w:r>'
$new27 = '<w:r><w:t xml:space="preserve">Update </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ApplicationDbContext</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IdentityDbContext</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> parent to use User class as </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IdentityUser</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> extender</w:t></w:r>'
if ($xml.IndexOf($old27) -lt 0) { throw 'pattern not found for para 27' }
$xml = $xml.Replace($old27, $new27)

$old30 = '<w:r><w:t>Add Identity Configuration to ServiceExtensions class</w:t></w:r>'
$new30 = '<w:r><w:t xml:space="preserve">Add Identity Configuration to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ServiceExtensions</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class</w:t></w:r>'
if ($xml.IndexOf($old30) -lt 0) { throw 'pattern not found for para 30' }
$xml = $xml.Replace($old30, $new30)

$old31 = '<w:r><w:t>Register IdentityConfiguration to the ConfigureServices method in the Startup class</w:t></w:r>'
$new31 = '<w:r><w:t xml:space="preserve">Register </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IdentityConfiguration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ConfigureServices</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method in the Startup class</w:t></w:r>'
if ($xml.IndexOf($old31) -lt 0) { throw 'pattern not found for para 31' }
$xml = $xml.Replace($old31, $new31)

$old32 = '<w:r><w:t>Register Authentication to ConfigureServices method in the Startup class</w:t></w:r>'
$new32 = '<w:r><w:t xml:space="preserve">Register Authentication to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ConfigureServices</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method in the Startup class</w:t></w:r>'
if ($xml.IndexOf($old32) -lt 0) { throw 'pattern not found for para 32' }
$xml = $xml.Replace($old32, $new32)

$old33 = '<w:r><w:t>Add UseAuthentication to the application’s request pipeline</w:t></w:r>'
$new33 = '<w:r><w:t xml:space="preserve">Add </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>UseAuthentication</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to the application’s request pipeline</w:t></w:r>'
if ($xml.IndexOf($old33) -lt 0) { throw 'pattern not found for para 33' }
$xml = $xml.Replace($old33, $new33)

$old36 = '<w:r><w:t>Add RolesConfiguration class to configuration folder</w:t></w:r>'
$new36 = '<w:r><w:t xml:space="preserve">Add </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RolesConfiguration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class to configuration folder</w:t></w:r>'
if ($xml.IndexOf($old36) -lt 0) { throw 'pattern not found for para 36' }
$xml = $xml.Replace($old36, $new36)

$old37 = '<w:r><w:t>Apply RolesConfiguration to the ApplicationDbContext’s OnModelCreating method</w:t></w:r>'
$new37 = '<w:r><w:t xml:space="preserve">Apply </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RolesConfiguration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ApplicationDbContext’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>OnModelCreating</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method</w:t></w:r>'
if ($xml.IndexOf($old37) -lt 0) { throw 'pattern not found for para 37' }
$xml = $xml.Replace($old37, $new37)

$oldTail = '<w:p w14:paraId="398886DF" w14:textId="25DD607A" w:rsidR="001A0698" w:rsidRDefault="001A0698" w:rsidP="00A55FB7"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Commit</w:t></w:r></w:p>'
$newTail = '<w:p w14:paraId="398886DF" w14:textId="25DD607A" w:rsidR="001A0698" w:rsidRDefault="001A0698" w:rsidP="00A55FB7"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Commit</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Add </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>JWTSettings</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>appsettings.json</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Create JWT configuration in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ServiceExtensions</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Register </w:t></w:r><w:r><w:t>JWT configuration</w:t></w:r><w:r><w:t xml:space="preserve"> to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ConfigureServices</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method in the Startup class</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Commit</w:t></w:r></w:p>'
if ($xml.IndexOf($oldTail) -lt 0) { throw 'pattern not found for tail' }
$xml = $xml.Replace($oldTail, $newTail)

$content.InsertXML($xml)

Write-Output ("Paragraphs:" + $d.Paragraphs.Count)
